# Regenerate save_data to use K (strikeouts) instead of Strike# for the
# "tepera_ryan" workbook. This updates column G (header "K") with the
# recalculated strikeout counts for each outing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new K value (column G)
$kValues = @{
    2  = 1
    3  = 0
    4  = 1
    5  = 2
    6  = 3
    7  = 2
    8  = 1
    9  = 1
    10 = 0
    11 = 0
    12 = 0
    14 = 1
    15 = 3
    16 = 2
    19 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
